$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift current rows 25-29 data down: the new row25 data is fresh,
# rows 26-29 take over what used to be in rows 25-28, and two new rows
# (30, 31) are appended carrying what used to be in rows 28 and 29.

# Row 25 (new values)
$ws.Range("D25").Value = 44663
$ws.Range("L25").Value = "Especial"
$ws.Range("M25").Value = 20
$ws.Range("N25").Value = 335000
$ws.Range("O25").Value = 340000
$ws.Range("P25").Value = 337500
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 750

# Row 26
$ws.Range("D26").Value = 44663
$ws.Range("N26").Value = 305000
$ws.Range("O26").Value = 310000
$ws.Range("P26").Value = 307500
$ws.Range("S26").Value = 683

# Row 27
$ws.Range("D27").Value = 44634
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 10
$ws.Range("R27").Value = "Región Metropolitana"

# Row 28
$ws.Range("D28").Value = 44622
$ws.Range("M28").Value = 16
$ws.Range("N28").Value = 410000
$ws.Range("O28").Value = 420000
$ws.Range("P28").Value = 415000
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 922

# Row 29
$ws.Range("L29").Value = "Especial"
$ws.Range("N29").Value = 305000
$ws.Range("O29").Value = 310000
$ws.Range("P29").Value = 307500
$ws.Range("R29").Value = "Provincia de Cachapoal"
$ws.Range("S29").Value = 683

# New row 30
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44309
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104003
$ws.Range("J30").Value = "Membrillo"
$ws.Range("K30").Value = "Champion"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 285000
$ws.Range("O30").Value = 290000
$ws.Range("P30").Value = 287500
$ws.Range("Q30").Value = "$/bins (450 kilos)"
$ws.Range("R30").Value = "Provincia de Cachapoal"
$ws.Range("S30").Value = 639
$ws.Range("T30").Value = 450

# New row 31
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44309
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100104
$ws.Range("H31").Value = "Frutos de pepita"
$ws.Range("I31").Value = 100104003
$ws.Range("J31").Value = "Membrillo"
$ws.Range("K31").Value = "Champion"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 20
$ws.Range("N31").Value = 255000
$ws.Range("O31").Value = 260000
$ws.Range("P31").Value = 257500
$ws.Range("Q31").Value = "$/bins (450 kilos)"
$ws.Range("R31").Value = "Provincia de Cachapoal"
$ws.Range("S31").Value = 572
$ws.Range("T31").Value = 450

# Apply the date format (same style as other date cells in column D)
$ws.Range("D30").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("D31").NumberFormat = $ws.Range("D29").NumberFormat
